# Added CO2 emission for cruises
# - adds two new worksheets that break out harbour -> country pairs
#   (one for "leg1_participants", one "_fr" variant with Las Palmas/Gran
#   Canaria as the first row instead of Stanley/Falkland Islands)
# - leaves Sheet1 as-is content-wise, just moves the active selection

$wb = $excel.ActiveWorkbook

# --- Sheet1: no data changes, just move the selection / drop tab focus ---
$ws1 = $wb.Worksheets.Item(1)
[void]$ws1.Range("E2").Select()

# --- New sheet 2: "2019_PS120_leg1_participants" ---
$ws2 = $wb.Worksheets.Add($null, $wb.Worksheets.Item($wb.Worksheets.Count))
$ws2.Name = "2019_PS120_leg1_participants"

$ws2.Range("A1").Value = "Stanley"
$ws2.Range("B1").Value = "Falkland Islands"
$ws2.Range("A2").Value = "Bremerhaven"
$ws2.Range("B2").Value = "Germany"
$ws2.Range("A3").Value = "Bremerhaven"
$ws2.Range("B3").Value = "Germany"
$ws2.Range("A4").Value = "Bremerhaven"
$ws2.Range("B4").Value = "Germany"
$ws2.Range("A5").Value = "Zagreb"
$ws2.Range("B5").Value = "Croatia"

$ws2.PageSetup.PaperSize = 9
$ws2.PageSetup.Orientation = 1

[void]$ws2.Range("A2:B5").Select()

# --- New sheet 3: "2019_PS120_leg1_participants_fr" ---
$ws3 = $wb.Worksheets.Add($null, $wb.Worksheets.Item($wb.Worksheets.Count))
$ws3.Name = "2019_PS120_leg1_participants_fr"

$ws3.Range("A1").Value = "Las Palmas"
$ws3.Range("B1").Value = "Gran Canaria"
$ws3.Range("A2").Value = "Bremerhaven"
$ws3.Range("B2").Value = "Germany"
$ws3.Range("A3").Value = "Bremerhaven"
$ws3.Range("B3").Value = "Germany"
$ws3.Range("A4").Value = "Bremerhaven"
$ws3.Range("B4").Value = "Germany"
$ws3.Range("A5").Value = "Zagreb"
$ws3.Range("B5").Value = "Croatia"

[void]$ws3.Range("C7").Select()

# The _fr sheet ends up the active tab
[void]$ws3.Activate()
